# Weekly fruit/vegetable price update:
# Insert a new weekly record (row 110) for "Camote" at "Vega Modelo de Temuco",
# pushing the existing rows 110-135 down to 111-136.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 110, shifting rows 110:135 down to 111:136.
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with the new weekly record.
$ws.Cells.Item(110, 1).Value = 10
$ws.Cells.Item(110, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(110, 3).Value = "La Araucanía"
$ws.Cells.Item(110, 4).Value = 44841
$ws.Cells.Item(110, 5).Value = 9
$ws.Cells.Item(110, 6).Value = 100114002
$ws.Cells.Item(110, 7).Value = "Camote"
$ws.Cells.Item(110, 8).Value = "Sin especificar"
$ws.Cells.Item(110, 9).Value = "Primera"
$ws.Cells.Item(110, 10).Value = 10
$ws.Cells.Item(110, 11).Value = 20000
$ws.Cells.Item(110, 12).Value = 20000
$ws.Cells.Item(110, 13).Value = 20000
$ws.Cells.Item(110, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(110, 15).Value = "Perú"
$ws.Cells.Item(110, 16).Value = 1000
$ws.Cells.Item(110, 17).Value = 20
$ws.Cells.Item(110, 18).Value = "Hortaliza"

# Apply the same date number format (as used by the rest of column D) to the new D110.
$ws.Cells.Item(110, 4).NumberFormat = $ws.Cells.Item(111, 4).NumberFormat
